# [IMP] new test data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = "####-<2-99"
$ws.Range("G14").Value = "####-<#-01"
$ws.Range("I14").Value = "####-<#-01"

$ws.Range("F15").Value = "####-<#-10"
$ws.Range("G15").Value = "####-<#-10"

$ws.Range("G17").Value = "####-<#-15"

$ws.Range("G18").Value = "####-<#-20"

$ws.Range("G19").Value = "####-<#-20"

$ws.Range("F20").Value = "####-<#-99"

$ws.Range("D12").Select()
